$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Förändrad" column (C) holds a date serial that was bulk-updated
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row.
$oldValue = 45172
$newValue = 45175

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 407
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
